# Commit: "added timings to SPI0 and SD_HST"
#
# Appends a trace-length-matching note to the Description (column J) of the
# SD_HST- and SPI0-related pin rows in the IO table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Description (J) gets " / +/-100ps inter-SD_HST" appended.
$sdHstRows = @(34, 38, 40, 42, 61, 63, 65)

# Rows whose Description (J) gets " / +/-100ps inter-SPI0" appended.
$spi0Rows = @(62, 64, 66, 68, 70, 72, 76)

foreach ($r in $sdHstRows) {
    $cell = $ws.Cells.Item($r, 10)
    $current = [string]$cell.Value2
    $cell.Value = $current.TrimEnd() + " / +/-100ps inter-SD_HST"
}

foreach ($r in $spi0Rows) {
    $cell = $ws.Cells.Item($r, 10)
    $current = [string]$cell.Value2
    $cell.Value = $current.TrimEnd() + " / +/-100ps inter-SPI0"
}
